# Scheduled runner update: refresh market-price derived columns (H-N) for
# affected Leve rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1477.25
$ws.Range("J112").Value = 1618.2858
$ws.Range("L112").Value = 4854.857400000001
$ws.Range("N112").Value = -7070.857400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3216.861
$ws.Range("I132").Value = 3318.3794
$ws.Range("J132").Value = 2796.2856
$ws.Range("K132").Value = 9955.138199999999
$ws.Range("L132").Value = 8388.856800000001
$ws.Range("M132").Value = -7425.138199999999
$ws.Range("N132").Value = -13448.8568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 70802.5
$ws.Range("I86").Value = 102026.73
$ws.Range("J86").Value = 2109.2
$ws.Range("K86").Value = 102026.73
$ws.Range("L86").Value = 2109.2
$ws.Range("M86").Value = -100903.73
$ws.Range("N86").Value = -4355.2

$ws.Range("H89").Value = 70802.5
$ws.Range("I89").Value = 102026.73
$ws.Range("J89").Value = 2109.2
$ws.Range("K89").Value = 510133.65
$ws.Range("L89").Value = 10546
$ws.Range("M89").Value = -504517.65
$ws.Range("N89").Value = -21778

$ws.Range("H107").Value = 83334360
$ws.Range("I107").Value = 200000480
$ws.Range("J107").Value = 1417.1428
$ws.Range("K107").Value = 200000480
$ws.Range("L107").Value = 1417.1428
$ws.Range("M107").Value = -199998560
$ws.Range("N107").Value = -5257.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 900.4
$ws.Range("I16").Value = 833.3333
$ws.Range("J16").Value = 1001
$ws.Range("K16").Value = 833.3333
$ws.Range("L16").Value = 1001
$ws.Range("M16").Value = -546.3333
$ws.Range("N16").Value = -1575

$ws.Range("H113").Value = 900.4
$ws.Range("I113").Value = 833.3333
$ws.Range("J113").Value = 1001
$ws.Range("K113").Value = 833.3333
$ws.Range("L113").Value = 1001
$ws.Range("M113").Value = 1336.6667
$ws.Range("N113").Value = -5341

$ws.Range("H132").Value = 4099.7334
$ws.Range("I132").Value = 3681.6365
$ws.Range("J132").Value = 5249.5
$ws.Range("K132").Value = 11044.9095
$ws.Range("L132").Value = 15748.5
$ws.Range("M132").Value = -8514.9095
$ws.Range("N132").Value = -20808.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1150
$ws.Range("I68").Value = 866.6667
$ws.Range("J68").Value = 1320
$ws.Range("K68").Value = 2600.0001
$ws.Range("L68").Value = 3960
$ws.Range("M68").Value = -1789.0001
$ws.Range("N68").Value = -5582

$ws.Range("H71").Value = 1150
$ws.Range("I71").Value = 866.6667
$ws.Range("J71").Value = 1320
$ws.Range("K71").Value = 7800.0003
$ws.Range("L71").Value = 11880
$ws.Range("M71").Value = -3744.0003
$ws.Range("N71").Value = -19992

$ws.Range("H75").Value = 1289.375
$ws.Range("I75").Value = 1000
$ws.Range("J75").Value = 1385.8334
$ws.Range("K75").Value = 3000
$ws.Range("L75").Value = 4157.5002
$ws.Range("M75").Value = -2002
$ws.Range("N75").Value = -6153.5002

$ws.Range("H78").Value = 1289.375
$ws.Range("I78").Value = 1000
$ws.Range("J78").Value = 1385.8334
$ws.Range("K78").Value = 9000
$ws.Range("L78").Value = 12472.5006
$ws.Range("M78").Value = -4008
$ws.Range("N78").Value = -22456.5006

$ws.Range("H92").Value = 750
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 3000
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -5496

$ws.Range("H107").Value = 376316
$ws.Range("I107").Value = 633.3333
$ws.Range("J107").Value = 869399.5
$ws.Range("K107").Value = 1899.9999
$ws.Range("L107").Value = 2608198.5
$ws.Range("M107").Value = 20.00009999999997
$ws.Range("N107").Value = -2612038.5

$ws.Range("H122").Value = 565
$ws.Range("J122").Value = 568.6667
$ws.Range("L122").Value = 5118.0003
$ws.Range("N122").Value = -10018.0003

$ws.Range("H131").Value = 816.86
$ws.Range("J131").Value = 822.0808
$ws.Range("L131").Value = 2466.2424
$ws.Range("N131").Value = -12546.2424

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 90911350
$ws.Range("I97").Value = 111113560
$ws.Range("J97").Value = 1400
$ws.Range("K97").Value = 111113560
$ws.Range("L97").Value = 1400
$ws.Range("M97").Value = -111113064
$ws.Range("N97").Value = -2392

$ws.Range("H102").Value = 275278.1
$ws.Range("I102").Value = 1788.1818
$ws.Range("J102").Value = 548768
$ws.Range("K102").Value = 1788.1818
$ws.Range("L102").Value = 548768
$ws.Range("M102").Value = -166.1818000000001
$ws.Range("N102").Value = -552012

$ws.Range("H113").Value = 1776.5555
$ws.Range("I113").Value = 1525
$ws.Range("J113").Value = 1848.4286
$ws.Range("K113").Value = 1525
$ws.Range("L113").Value = 1848.4286
$ws.Range("M113").Value = 645
$ws.Range("N113").Value = -6188.4286

$ws.Range("H126").Value = 3269976
$ws.Range("I126").Value = 2406.375
$ws.Range("J126").Value = 5884031.5
$ws.Range("K126").Value = 7219.125
$ws.Range("L126").Value = 17652094.5
$ws.Range("M126").Value = -4749.125
$ws.Range("N126").Value = -17657034.5

$ws.Range("H132").Value = 2338.1333
$ws.Range("I132").Value = 2130.5417
$ws.Range("J132").Value = 3168.5
$ws.Range("K132").Value = 6391.625100000001
$ws.Range("L132").Value = 9505.5
$ws.Range("M132").Value = -3861.625100000001
$ws.Range("N132").Value = -14565.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4467.6665
$ws.Range("I7").Value = 2800.8
$ws.Range("J7").Value = 6551.25
$ws.Range("K7").Value = 2800.8
$ws.Range("L7").Value = 6551.25
$ws.Range("M7").Value = -2688.8
$ws.Range("N7").Value = -6775.25

$ws.Range("H122").Value = 2319.05
$ws.Range("I122").Value = 2243.389
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 6730.167
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -4280.167
$ws.Range("N122").Value = -13900

$ws.Range("H126").Value = 4467.6665
$ws.Range("I126").Value = 2800.8
$ws.Range("J126").Value = 6551.25
$ws.Range("K126").Value = 8402.400000000001
$ws.Range("L126").Value = 19653.75
$ws.Range("M126").Value = -5932.400000000001
$ws.Range("N126").Value = -24593.75

$ws.Range("H132").Value = 4952.2856
$ws.Range("I132").Value = 4764.706
$ws.Range("K132").Value = 14294.118
$ws.Range("M132").Value = -11764.118

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H136").Value = 1470.3235
$ws.Range("J136").Value = 10668
$ws.Range("L136").Value = 32004
$ws.Range("N136").Value = -37104
